$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Update Price (D) and Volume(1h) (E) columns for rows with changed values
Set-TextValue $ws.Range("D2") "28.981.44"
$ws.Range("E2").Value = "  -0.22%  "
Set-TextValue $ws.Range("D3") "1.825.77"
$ws.Range("E3").Value = "  -0.28%  "
Set-TextValue $ws.Range("D4") "0.9962"
$ws.Range("E4").Value = "  -0.26%  "
Set-TextValue $ws.Range("D5") "243.70"
$ws.Range("E5").Value = "  +0.89%  "
Set-TextValue $ws.Range("D6") "0.6312"
$ws.Range("E6").Value = "  +0.70%  "
Set-TextValue $ws.Range("D7") "0.9984"
$ws.Range("E7").Value = "  -0.18%  "
Set-TextValue $ws.Range("D8") "0.07512"
$ws.Range("E8").Value = "  -1.31%  "
Set-TextValue $ws.Range("D9") "0.2940"
$ws.Range("E9").Value = "  +0.78%  "
Set-TextValue $ws.Range("D10") "23.05"
$ws.Range("E10").Value = "  +1.06%  "
Set-TextValue $ws.Range("D11") "0.07695"
$ws.Range("E11").Value = "  +0.72%  "
Set-TextValue $ws.Range("D12") "1.827.58"
$ws.Range("E12").Value = "  -0.18%  "
Set-TextValue $ws.Range("D13") "4.990"
$ws.Range("E13").Value = "  +0.63%  "
Set-TextValue $ws.Range("D14") "0.6685"
$ws.Range("E14").Value = "  +0.49%  "
Set-TextValue $ws.Range("D15") "83.01"
$ws.Range("E15").Value = "  +0.77%  "
Set-TextValue $ws.Range("D16") "0.000009626"
$ws.Range("E16").Value = "  +1.27%  "
Set-TextValue $ws.Range("D17") "6.045"
$ws.Range("E17").Value = "  +0.94%  "
Set-TextValue $ws.Range("D18") "29.026.26"
$ws.Range("E18").Value = "  +0.21%  "
Set-TextValue $ws.Range("D19") "12.57"
$ws.Range("E19").Value = "  +1.95%  "
Set-TextValue $ws.Range("D20") "226.07"
$ws.Range("E20").Value = "  +0.46%  "
Set-TextValue $ws.Range("D21") "0.9977"
$ws.Range("E21").Value = "  -0.17%  "
Set-TextValue $ws.Range("D22") "7.146"
$ws.Range("E22").Value = "  -1.03%  "
Set-TextValue $ws.Range("D23") "0.9978"
$ws.Range("E23").Value = "  -0.28%  "
Set-TextValue $ws.Range("D24") "160.38"
$ws.Range("E24").Value = "  -0.46%  "
Set-TextValue $ws.Range("D25") "0.1425"
$ws.Range("E25").Value = "  +4.53%  "
Set-TextValue $ws.Range("D26") "8.501"
$ws.Range("E26").Value = "  +0.98%  "
Set-TextValue $ws.Range("D28") "1.499"
$ws.Range("E28").Value = "  +0.23%  "
Set-TextValue $ws.Range("D29") "4.138"
$ws.Range("E29").Value = "  +2.03%  "
Set-TextValue $ws.Range("D30") "4.058"
$ws.Range("E30").Value = "  +0.55%  "
Set-TextValue $ws.Range("D31") "0.05467"
$ws.Range("E31").Value = "  +5.14%  "
Set-TextValue $ws.Range("D32") "1.201"
$ws.Range("E32").Value = "  +0.44%  "
Set-TextValue $ws.Range("D33") "1.854"
$ws.Range("E33").Value = "  +0.22%  "
Set-TextValue $ws.Range("D34") "0.7439"
$ws.Range("E34").Value = "  +2.07%  "
Set-TextValue $ws.Range("D36") "2.644"
$ws.Range("E36").Value = "  +1.57%  "
Set-TextValue $ws.Range("D37") "1.240.85"
$ws.Range("E37").Value = "  -2.58%  "
Set-TextValue $ws.Range("D38") "2.748"
$ws.Range("E38").Value = "  -0.41%  "
Set-TextValue $ws.Range("D39") "0.01777"
$ws.Range("E39").Value = "  -0.46%  "
Set-TextValue $ws.Range("D40") "6.652"
$ws.Range("E40").Value = "  +2.19%  "
Set-TextValue $ws.Range("D41") "0.8985"
$ws.Range("E41").Value = "  +0.89%  "
Set-TextValue $ws.Range("D42") "0.9987"
$ws.Range("E42").Value = "  -0.15%  "
Set-TextValue $ws.Range("D43") "101.29"
$ws.Range("E43").Value = "  -0.11%  "
Set-TextValue $ws.Range("D46") "65.21"
$ws.Range("E46").Value = "  +2.34%  "
Set-TextValue $ws.Range("D47") "0.5083"
$ws.Range("E47").Value = "  -0.44%  "
Set-TextValue $ws.Range("D48") "0.4053"
$ws.Range("E48").Value = "  +1.85%  "
Set-TextValue $ws.Range("D49") "8.948"
$ws.Range("E49").Value = "  +1.74%  "
Set-TextValue $ws.Range("D50") "1.651"
$ws.Range("E50").Value = "  +0.72%  "
Set-TextValue $ws.Range("D51") "0.05782"
$ws.Range("E51").Value = "  +0.35%  "

# Rows where only Volume(1h) changed (Price unchanged)
$ws.Range("E27").Value = "  +0.43%  "
$ws.Range("E35").Value = "  -1.34%  "

# Rows 44/45: RocketPoolETH and BabyDogeCoin swap positions, with updated values
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextValue $ws.Range("D44") "1.981.04"
$ws.Range("E44").Value = "  +0.30%  "

$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws.Range("D45") "0.00000000123"
$ws.Range("E45").Value = "  +2.96%  "
